$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Turn on the AutoFilter over E1:E100 FIRST (while the sheet still
#    only has data through row 61) so the persisted <autoFilter> ref
#    stays pinned at E1:E100 instead of auto-growing to the new extent
#    once column E gets populated down to row 132 below.
# ---------------------------------------------------------------------
$ws.Range("E1:E100").AutoFilter()

$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Tabelle1!`$E`$1:`$E`$100")
$filterName.Visible = $false

# ---------------------------------------------------------------------
# 2) New data rows (57-132): completes model #6 and adds models #7-#14.
#    Each tuple: (row, A model#, B model name, C trloss, D val loss,
#                 E epoch, F saved(y/n), G metric label, H metric value)
# ---------------------------------------------------------------------
$rows = @(
    @(57, 6, "dense128_relu_e-5_dropout03_dense64_relu_e-7", 441, 493, 3, "y", $null, $null),
    @(58, 6, "dense128_relu_e-5_dropout03_dense64_relu_e-7", 377, 501, 4, "n", $null, $null),
    @(59, 6, "dense128_relu_e-5_dropout03_dense64_relu_e-7", 334, 490, 5, "y", $null, $null),
    @(60, 6, "dense128_relu_e-5_dropout03_dense64_relu_e-7", 309, 528, 6, "n", $null, $null),
    @(61, 6, "dense128_relu_e-5_dropout03_dense64_relu_e-7", 275, 514, 7, "n", $null, $null),
    @(62, 6, "dense128_relu_e-5_dropout03_dense64_relu_e-7", 260, 546, 8, "n", "teacc", 1437),
    @(63, 7, "dense128_relu_e-5_dropout03_dense64_relu_e-5", 705, 641, 1, "y", $null, $null),
    @(64, 7, "dense128_relu_e-5_dropout03_dense64_relu_e-5", 533, 585, 2, "y", $null, $null),
    @(65, 7, "dense128_relu_e-5_dropout03_dense64_relu_e-5", 443, 510, 3, "y", $null, $null),
    @(66, 7, "dense128_relu_e-5_dropout03_dense64_relu_e-5", 386, 476, 4, "y", $null, $null),
    @(67, 7, "dense128_relu_e-5_dropout03_dense64_relu_e-5", 345, 485, 5, "n", $null, $null),
    @(68, 7, "dense128_relu_e-5_dropout03_dense64_relu_e-5", 311, 516, 6, "n", "teacc", 6056),
    @(69, 7, "dense128_relu_e-5_dropout03_dense64_relu_e-5", 290, 517, 7, "n", "roc", 1427),
    @(70, 8, "dense128_relu_dense128_relu", 637, 493, 1, "y", $null, $null),
    @(71, 8, "dense128_relu_dense128_relu", 302, 490, 2, "y", $null, $null),
    @(72, 8, "dense128_relu_dense128_relu", 175, 459, 3, "y", $null, $null),
    @(73, 8, "dense128_relu_dense128_relu", 114, 520, 4, "n", "teacc", 7073),
    @(74, 8, "dense128_relu_dense128_relu", 83, 516, 5, "n", "roc", 1297),
    @(75, 9, "dense128_relu_dense128_relu_dense128_relu", 679, 502, 1, "y", $null, $null),
    @(76, 9, "dense128_relu_dense128_relu_dense128_relu", 318, 466, 2, "y", $null, $null),
    @(77, 9, "dense128_relu_dense128_relu_dense128_relu", 195, 465, 3, "y", $null, $null),
    @(78, 9, "dense128_relu_dense128_relu_dense128_relu", 120, 536, 4, "n", "teacc", 6627),
    @(79, 9, "dense128_relu_dense128_relu_dense128_relu", 89, 598, 5, "n", "roc", 1422),
    @(80, 10, "d128_relu_dr1_d128_relu_dr1_d128", 644, 568, 1, "y", $null, $null),
    @(81, 10, "d128_relu_dr1_d128_relu_dr1_d128", 427, 516, 2, "y", $null, $null),
    @(82, 10, "d128_relu_dr1_d128_relu_dr1_d128", 306, 491, 3, "y", $null, $null),
    @(83, 10, "d128_relu_dr1_d128_relu_dr1_d128", 218, 480, 4, "y", $null, $null),
    @(84, 10, "d128_relu_dr1_d128_relu_dr1_d128", 153, 545, 5, "n", "teacc", 6459),
    @(85, 10, "d128_relu_dr1_d128_relu_dr1_d128", 113, 616, 6, "n", "roc", 1412),
    @(86, 11, "d128_relu_e-5_dr2_d128_relu_dr2_d128", 704, 810, 1, "y", $null, $null),
    @(87, 11, "d128_relu_e-5_dr2_d128_relu_dr2_d128", 525, 639, 2, "y", $null, $null),
    @(88, 11, "d128_relu_e-5_dr2_d128_relu_dr2_d128", 455, 499, 3, "y", $null, $null),
    @(89, 11, "d128_relu_e-5_dr2_d128_relu_dr2_d128", 390, 582, 4, "n", $null, $null),
    @(90, 11, "d128_relu_e-5_dr2_d128_relu_dr2_d128", 344, 550, 5, "n", $null, $null),
    @(91, 11, "d128_relu_e-5_dr2_d128_relu_dr2_d128", 300, 654, 6, "n", "teacc", 5459),
    @(92, 11, "d128_relu_e-5_dr2_d128_relu_dr2_d128", 269, 649, 7, "n", "roc", 1543),
    @(93, 12, "dense128_relu_5e-5", 775, 508, 1, "y", $null, $null),
    @(94, 12, "dense128_relu_5e-5", 441, 466, 2, "y", $null, $null),
    @(95, 12, "dense128_relu_5e-5", 343, 437, 3, "y", $null, $null),
    @(96, 12, "dense128_relu_5e-5", 295, 415, 4, "y", $null, $null),
    @(97, 12, "dense128_relu_5e-5", 275, 415, 5, "n", $null, $null),
    @(98, 12, "dense128_relu_5e-5", 252, 417, 6, "y", $null, $null),
    @(99, 12, "dense128_relu_5e-5", 239, 409, 7, "y", $null, $null),
    @(100, 12, "dense128_relu_5e-5", 232, 407, 8, "y", $null, $null),
    @(101, 12, "dense128_relu_5e-5", 228, 416, 9, "n", $null, $null),
    @(102, 12, "dense128_relu_5e-5", 223, 402, 10, "y", $null, $null),
    @(103, 12, "dense128_relu_5e-5", 222, 400, 11, "y", $null, $null),
    @(104, 12, "dense128_relu_5e-5", 218, 408, 12, "n", $null, $null),
    @(105, 12, "dense128_relu_5e-5", 214, 401, 13, "y", $null, $null),
    @(106, 12, "dense128_relu_5e-5", 210, 413, 14, "n", $null, $null),
    @(107, 12, "dense128_relu_5e-5", 209, 398, 15, "y", $null, $null),
    @(108, 12, "dense128_relu_5e-5", 205, 400, 16, "n", "teacc", 9065),
    @(109, 12, "dense128_relu_5e-5", 203, 403, 17, "n", "roc", 1264),
    @(110, 13, "dense64_relu", 664, 476, 1, "y", $null, $null),
    @(111, 13, "dense64_relu", 273, 452, 2, "y", $null, $null),
    @(112, 13, "dense64_relu", 158, 448, 3, "y", $null, $null),
    @(113, 13, "dense64_relu", 107, 438, 4, "y", $null, $null),
    @(114, 13, "dense64_relu", 78, 441, 5, "n", $null, $null),
    @(115, 13, "dense64_relu", 62, 460, 6, "n", "teacc", 8679),
    @(116, 13, "dense64_relu", 52, 482, 7, "n", "roc", 1371),
    @(117, 14, "dense64_relu_e-5", 687, 461, 1, "y", $null, $null),
    @(118, 14, "dense64_relu_e-5", 313, 434, 2, "y", $null, $null),
    @(119, 14, "dense64_relu_e-5", 212, 422, 3, "y", $null, $null),
    @(120, 14, "dense64_relu_e-5", 168, 432, 4, "n", $null, $null),
    @(121, 14, "dense64_relu_e-5", 148, 412, 5, "y", $null, $null),
    @(122, 14, "dense64_relu_e-5", 131, 406, 6, "y", $null, $null),
    @(123, 14, "dense64_relu_e-5", 122, 418, 7, "n", $null, $null),
    @(124, 14, "dense64_relu_e-5", 113, 412, 8, "n", $null, $null),
    @(125, 14, "dense64_relu_e-5", 110, 411, 9, "n", $null, $null),
    @(126, 14, "dense64_relu_e-5", 103, 410, 10, "n", $null, $null),
    @(127, 14, "dense64_relu_e-5", 98, 408, 11, "y", $null, $null),
    @(128, 14, "dense64_relu_e-5", 95, 416, 12, "n", $null, $null),
    @(129, 14, "dense64_relu_e-5", 92, 407, 13, "y", $null, $null),
    @(130, 14, "dense64_relu_e-5", 90, 407, 14, "n", $null, $null),
    @(131, 14, "dense64_relu_e-5", 89, 407, 15, "n", "teacc", 9362),
    @(132, 14, "dense64_relu_e-5", 87, 412, 16, "n", "roc", 1258)
)

foreach ($item in $rows) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
    $ws.Cells.Item($r, 5).Value = $item[5]
    if ($item[6] -ne $null) {
        $ws.Cells.Item($r, 6).Value = $item[6]
    }
    if ($item[7] -ne $null) {
        $ws.Cells.Item($r, 7).Value = $item[7]
    }
    if ($item[8] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $item[8]
    }
}

# ---------------------------------------------------------------------
# 3) Update the view: scroll/select near the new bottom of the data and
#    select the next empty cell in column I (mirrors the saved state).
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A108"))
$ws.Range("I134").Select()
